$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: row, column, expected old text, new text.
# Using direct table-cell addressing (rather than a global Find/Replace)
# avoids ambiguity because several expressions share the same "old" text
# (e.g. "56÷6=" appears twice) while having different replacements.
$replacements = @(
    @(1, 1, "37÷8=", "56÷6="),
    @(1, 2, "31÷2=", "82÷9="),
    @(1, 3, "86÷5=", "99÷5="),
    @(1, 4, "29÷9=", "53÷6="),
    @(1, 5, "22÷4=", "43÷8="),
    @(5, 1, "85÷9=", "63÷9="),
    @(5, 2, "58÷9=", "70÷4="),
    @(5, 3, "14÷9=", "25÷9="),
    @(5, 4, "74÷7=", "56÷8="),
    @(5, 5, "96÷9=", "85÷6="),
    @(9, 1, "26÷3=", "57÷7="),
    @(9, 2, "56÷6=", "88÷4="),
    @(9, 3, "64÷9=", "78÷5="),
    @(9, 4, "21÷5=", "90÷3="),
    @(9, 5, "59÷5=", "73÷7="),
    @(13, 1, "78÷3=", "48÷2="),
    @(13, 2, "63÷8=", "87÷5="),
    @(13, 3, "30÷6=", "56÷2="),
    @(13, 4, "68÷7=", "38÷3="),
    @(13, 5, "80÷3=", "90÷6="),
    @(17, 1, "29÷4=", "98÷4="),
    @(17, 2, "78÷7=", "77÷9="),
    @(17, 3, "64÷4=", "16÷9="),
    @(17, 4, "27÷6=", "84÷6="),
    @(17, 5, "84÷7=", "77÷9=")
)

foreach ($item in $replacements) {
    $row = $item[0]
    $col = $item[1]
    $old = $item[2]
    $new = $item[3]

    $cell = $t.Cell($row, $col)
    $cellRange = $cell.Range
    $cellRange.End = $cellRange.End - 1

    if ($cellRange.Text -eq $old) {
        $cellRange.Text = $new
    } else {
        Write-Output "Unexpected text at row $row col $col`: $($cellRange.Text) (expected $old)"
    }
}
